$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the two backup codes that currently live further down the sheet
# (rows 20/21 - NHYK5008HQDA and ZADGNDVPP03M) before we clear them out.
$code11 = $ws.Range("A20").Value2
$code12 = $ws.Range("A21").Value2

# Remove the old/expired backup codes in A2:A4 and replace with the ones
# that used to sit in A17:A19.
$ws.Range("A2").Value = $ws.Range("A17").Value2
$ws.Range("A3").Value = $ws.Range("A18").Value2
$ws.Range("A4").Value = $ws.Range("A19").Value2

# Clear out the old rows 17-21 entirely.
$ws.Range("A17:A21").Clear()

# Put the saved codes into the new, shorter list (A11:A12).
$ws.Range("A11").Value = $code11
$ws.Range("A12").Value = $code12

# Update the selected cell to match the author's saved view state.
$ws.Range("C17").Select()
